# Insert a new weekly price-report row above the current row 165 for
# "Terminal La Palmera de La Serena" (Zapallo italiano), shifting the
# existing rows 165:188 down to 166:189, then populate the newly
# inserted row 165 with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 165 and below down by one row.
$ws.Rows.Item(165).Insert()

# Fill in the new row with the latest observation.
$ws.Cells.Item(165, 1).Value2  = 8
$ws.Cells.Item(165, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(165, 3).Value2  = "Coquimbo"
$ws.Cells.Item(165, 4).Value2  = 44474
$ws.Cells.Item(165, 5).Value2  = 4
$ws.Cells.Item(165, 6).Value2  = 100112032
$ws.Cells.Item(165, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(165, 8).Value2  = "Sin especificar"
$ws.Cells.Item(165, 9).Value2  = "Primera"
$ws.Cells.Item(165, 10).Value2 = 400
$ws.Cells.Item(165, 11).Value2 = 19000
$ws.Cells.Item(165, 12).Value2 = 20000
$ws.Cells.Item(165, 13).Value2 = 19500
$ws.Cells.Item(165, 14).Value2 = "$/caja 70 unidades"
$ws.Cells.Item(165, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(165, 16).Value2 = 279
$ws.Cells.Item(165, 17).Value2 = 70
$ws.Cells.Item(165, 18).Value2 = "Hortaliza"
